# Updated queries for C3DC first half testcases.
#
# The embedded SQL statements joined df_study/df_participant/df_diagnoses/
# df_treatments/df_treatment_resp/df_survival/df_reference_files on bare
# ".id" columns. The join keys were renamed to the fully-qualified
# "<table>_id" columns (std.study_id / prt.participant_id, and the
# matching quoted "table.column" aliases on the right-hand side).
#
# Apply the same set of textual substitutions to every cell on the sheet
# that contains one of the affected SELECT statements.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells B2, C2, B3, B4, B5, B6 and B7 hold the long SQL query strings that
# reference the old join columns.
$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $v = $rng.Value2

    if ($v -ne $null) {
        $v = $v.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
        $v = $v.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
        $v = $v.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
        $v = $v.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
        $v = $v.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
        $v = $v.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

        $rng.Value = $v
    }
}
